$d = $word.ActiveDocument

# --- Edit 1: "Technische interesses" paragraph about PHP/L-Systems project ---
$old1 = "Samen hebben we dit vak afgerond met een 10, omdat we zonder elkaar geen idee hadden voor een opdracht."
$new1 = "Samen hebben we dit vak afgerond met een 10. We hebben er samen aan gewerkt, omdat we zonder elkaar geen idee hadden voor een opdracht."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# --- Edit 2: Android telefoon paragraph ---
$old2 = "De huidige ontwikkelomgevingen zoals NetBeans en Eclipse kunnen gebruikt worden voor het programmeren van de applicaties. Verder is er een geïnstalleerde Android SDK nodig, en een Virtual Device."
$new2 = "Voor programmeurs is het Java platform één van de beste omgevingen om in te werken. Je ziet dat mobiele applicaties steeds belangrijker worden in de maatschappij. Daardoor stijgt de vraag naar kwalitatief goede mobiele software. Ik zet mij graag in om de eisen en wensen van de mobiele consument waar te maken."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
